# Rename the paired "_old"/"_new" column-header suffixes to the
# format-version-specific suffixes "_FV2310" / "_FV2404", freeze the
# header row, and turn the used range into a proper Excel Table - matching
# the "adapt column header formatting to respective input file names" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J (1..10) carried the "_old" suffix -> "_FV2310"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2310"
}

# Column K (11) is the untouched "diff" header.

# Columns L..U (12..21) carried the "_new" suffix -> "_FV2404"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2404"
}

# --- Turn the data range into an Excel Table ("Table1") -------------------
# The header row already carries hand-rolled direct formatting (bold, grey
# fill, borders, centered+wrapped). Converting the range straight into a
# ListObject would make Excel compute a header-row "dxf" to reconcile that
# existing formatting with the new table style, polluting styles.xml with
# an extra <dxfs> entry the authored workbook doesn't have. So: stash the
# header row's direct formatting, blank its named style to "Normal" (so the
# table gets created against a plain header), add the table, then restore
# the original direct formatting on top.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")

$headerRange.Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.Style = "Normal"

$dataRange = $ws.Range("A1:U60")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$scratch.Clear()

# --- Freeze the header row -------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
